$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 104905116.26
$ws.Range("P2").Value = 945745183.5700001
$ws.Range("Q2").Value = 818702785.8200001
$ws.Range("R2").Value = 6.6574671297
$ws.Range("S2").Value = 714440488.4299999
$ws.Range("T2").Value = 714440488.4299999
$ws.Range("U2").Value = 4.2751402151
$ws.Range("V2").Value = 34282999.37
$ws.Range("W2").Value = 24995071.91
$ws.Range("X2").Value = 1736159.66
$ws.Range("Y2").Value = 120292233.1
$ws.Range("Z2").Value = 119887919.82
$ws.Range("AA2").Value = 14982803.56
$ws.Range("AG2").Value = 5925016.79
$ws.Range("AP2").Value = 9.108964307899999
$ws.Range("AQ2").Value = 22.0233846542
$ws.Range("AR2").Value = 18.543399158102
$ws.Range("AS2").Value = 81090116.26000001
$ws.Range("AT2").Value = 3.265927006884
